$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Header text tweaks (rich-text shared strings) -- update only the digits
# that changed, leaving surrounding runs/formatting untouched.
# -----------------------------------------------------------------------

# "Volume 30   Number  14" -> "...15"  (A8, 1-based chars 21-22 = "14")
$ws.Range("A8").Characters(21, 2).Text = "15"

# "Report Covering the Week  4/3/2023  Through  4/9/2023"
#   -> "...4/10/2023  Through  4/16/2023"
# Edit the later date first so the earlier date's character offset doesn't shift.
$ws.Range("C9").Characters(46, 8).Text = "4/16/2023"
$ws.Range("C9").Characters(27, 8).Text = "4/10/2023"

# -----------------------------------------------------------------------
# Weekly crime-stat numbers (rows 15-30)
# -----------------------------------------------------------------------

# Row 15
$ws.Range("L15").Value = 0

# Row 16
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = -18.840579710144
$ws.Range("L16").Value = 19.148936170212
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = -77.865612648221

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 7.692307692307
$ws.Range("I17").Value = 49
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 22.5
$ws.Range("L17").Value = 58.064516129032
$ws.Range("M17").Value = 81.481481481481
$ws.Range("N17").Value = -16.949152542372

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = -26.470588235294
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = -16
$ws.Range("L18").Value = 47.368421052631
$ws.Range("M18").Value = 33.333333333333
$ws.Range("N18").Value = -60.563380281690

# Row 19
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -7.692307692307
$ws.Range("G19").Value = 101
$ws.Range("H19").Value = -16.831683168316
$ws.Range("I19").Value = 353
$ws.Range("J19").Value = 303
$ws.Range("K19").Value = 16.501650165016
$ws.Range("L19").Value = 143.448275862069
$ws.Range("M19").Value = 27.436823104693
$ws.Range("N19").Value = -44.929797191887

# Row 20
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = -25
$ws.Range("L20").Value = 200
$ws.Range("M20").Value = -10
$ws.Range("N20").Value = -96.103896103896

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -21.951219512195
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 173
$ws.Range("H21").Value = -17.919075144508
$ws.Range("I21").Value = 553
$ws.Range("J21").Value = 529
$ws.Range("K21").Value = 4.536862003780
$ws.Range("L21").Value = 94.035087719298
$ws.Range("M21").Value = 31.980906921241
$ws.Range("N21").Value = -60.556348074179

# Row 22: C22 changes from a plain number (1) to the text marker "0",
# matching the style already used by the other "N/A" cells in the row (D22).
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("L22").Value = 100

# Row 24
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = -31.578947368421
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 160
$ws.Range("H24").Value = -11.25
$ws.Range("I24").Value = 500
$ws.Range("J24").Value = 491
$ws.Range("K24").Value = 1.832993890020
$ws.Range("L24").Value = 58.227848101265
$ws.Range("M24").Value = 24.069478908188

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -11.111111111111
$ws.Range("F25").Value = 36
$ws.Range("H25").Value = -18.181818181818
$ws.Range("I25").Value = 121
$ws.Range("J25").Value = 95
$ws.Range("K25").Value = 27.368421052631
$ws.Range("L25").Value = 163.04347826087
$ws.Range("M25").Value = 86.153846153846

# Row 26
$ws.Range("L26").Value = 100

# Row 27: D27 (1 -> text "0") and E27 (0 -> text "***.*") gain the "N/A"
# marker style already used elsewhere on the sheet (row 23 cells).
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("E23").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 24
$ws.Range("K27").Value = 41.176470588235
$ws.Range("L27").Value = 200

# Row 30: G30 (1 -> text "0") and H30 (-100 -> text "***.*")
$ws.Range("G30").Value = "'0"
$ws.Range("F30").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = "***.*"
$ws.Range("E23").Copy()
$ws.Range("H30").PasteSpecial(-4122)
